$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.609.43"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "2.279.30"
$ws.Range("E3").Value = "  +3.05%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'251.35"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").Value = "'0.637"
$ws.Range("E6").Value = "  +2.40%  "
$ws.Range("D7").Value = "'73.41"
$ws.Range("E7").Value = "  +8.73%  "
$ws.Range("D9").Value = "'0.641"
$ws.Range("E9").Value = "  +3.23%  "
$ws.Range("D10").Value = "'39.06"
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").Value = "'0.0982"
$ws.Range("E11").Value = "  +4.84%  "
$ws.Range("D12").Value = "'59.00"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("E13").Value = "  +4.67%  "
$ws.Range("D14").Value = "'0.106"
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").Value = "2.618.79"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").Value = "'14.99"
$ws.Range("E16").Value = "  +3.31%  "
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "2.290.88"
$ws.Range("E18").Value = "  +4.37%  "
$ws.Range("D19").Value = "42.515.19"
$ws.Range("E19").Value = "  +1.43%  "
$ws.Range("E20").Value = "  +5.00%  "
$ws.Range("D21").Value = "'6.30"
$ws.Range("E21").Value = "  +2.56%  "
$ws.Range("D22").Value = "'72.19"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "'232.65"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").Value = "'2.20"
$ws.Range("E24").Value = "  +9.11%  "
$ws.Range("D25").Value = "'3.91"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").Value = "'11.51"
$ws.Range("E26").Value = "  +2.80%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "'3.65"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("D30").Value = "'2.14"
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").Value = "'166.57"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").Value = "'21.07"
$ws.Range("E32").Value = "  +3.29%  "
$ws.Range("E33").Value = "  +7.91%  "
$ws.Range("E34").Value = "  +4.95%  "
$ws.Range("D35").Value = "'0.0811"
$ws.Range("E35").Value = "  +3.61%  "
$ws.Range("D36").Value = "'31.62"
$ws.Range("E36").Value = "  +21.70%  "
$ws.Range("E37").Value = "  +3.58%  "
$ws.Range("D38").Value = "'4.67"
$ws.Range("E38").Value = "  +14.33%  "
$ws.Range("E39").Value = "  +3.40%  "
$ws.Range("D40").Value = "'0.0307"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'13.82"
$ws.Range("E41").Value = "  +15.42%  "
$ws.Range("E42").Value = "  +4.94%  "
$ws.Range("E43").Value = "  +4.94%  "
$ws.Range("D44").Value = "'0.212"
$ws.Range("E44").Value = "  +8.69%  "
$ws.Range("E45").Value = "  +6.76%  "
$ws.Range("D46").Value = "'61.74"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("E47").Value = "  -4.64%  "
$ws.Range("E48").Value = "  +4.43%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").Value = "'1.17"
$ws.Range("E50").Value = "  +1.91%  "
$ws.Range("D51").Value = "'97.70"
$ws.Range("E51").Value = "  +4.76%  "
